# Trade #7 closed at 2026-02-17 23:52:40 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.06   # Current Capital
$summary.Range("B4").Value = 0.06      # Total P&L $
$summary.Range("B5").Value = 0.17      # Total P&L %
$summary.Range("B6").Value = 7         # Total Trades
$summary.Range("B8").Value = 3         # Losing Trades
$summary.Range("B9").Value = 42.86     # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.06
$status.Range("D6").Value = 7
$status.Range("E6").Value = 0.06
$status.Range("F6").Value = 0.06
$status.Range("G6").Value = 42.86

# --- New trade row (#7) to append to "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 7
    B = "2026-02-17"
    C = "23:52:34"
    D = "MarketMaking"
    E = "UP"
    F = 0.9399999999999999
    G = 0.86
    H = "CLOSED"
    I = -8.5106
    J = -0.08
    K = 100.06
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A8").Value = $newRow.A
    # Prefix with an apostrophe so the date-formatted string is stored as
    # literal text (matching the source row's inline-string date cells)
    # instead of being auto-converted to a date serial number, then
    # restore the default "Normal" style so no stray number format sticks.
    $ws.Range("B8").Value = "'" + $newRow.B
    $ws.Range("B8").Style = "Normal"
    $ws.Range("C8").Value = $newRow.C
    $ws.Range("D8").Value = $newRow.D
    $ws.Range("E8").Value = $newRow.E
    $ws.Range("F8").Value = $newRow.F
    $ws.Range("G8").Value = $newRow.G
    $ws.Range("H8").Value = $newRow.H
    $ws.Range("I8").Value = $newRow.I
    $ws.Range("J8").Value = $newRow.J
    $ws.Range("K8").Value = $newRow.K
    $ws.Range("L8").Value = $newRow.L
    $ws.Range("M8").Value = $newRow.M
    $ws.Range("N8").Value = $newRow.N
    $ws.Range("O8").Value = $newRow.O
    $ws.Range("P8").Value = $newRow.P
    $ws.Range("Q8").Value = $newRow.Q
}
